$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): reorder column labels B..G
$ws.Range("B1").Value = "start_end"
$ws.Range("C1").Value = "b3_-"
$ws.Range("D1").Value = "4_-"
$ws.Range("E1").Value = "b6_"
$ws.Range("F1").Value = "b6_7"
$ws.Range("G1").Value = "7_"

# Row 2 (start_end)
$ws.Range("A2").Value = "start_end"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

# Row 3 (b3_-)
$ws.Range("A3").Value = "b3_-"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0

# Row 4 (4_-)
$ws.Range("A4").Value = "4_-"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1

# Row 5 (b6_)
$ws.Range("A5").Value = "b6_"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

# Row 6 (b6_7)
$ws.Range("A6").Value = "b6_7"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

# Row 7 (7_)
$ws.Range("A7").Value = "7_"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
